$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.434.38'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.477.34'
$ws.Range("E3").Value = '  -2.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").Value = '2.476.48'
$ws.Range("E9").Value = '  -2.73%  '
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.18'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").Value = '2.955.45'
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").Value = '67.429.28'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '2.481.83'
$ws.Range("E18").Value = '  -2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.55%  '
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '366.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.20%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.92%  '
$ws.Range("E26").Value = '  -4.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.32%  '
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").Value = '2.621.20'
$ws.Range("D30").Value = '0.0₃0957'
$ws.Range("E30").Value = '  -2.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '529.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -3.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.20%  '
$ws.Range("E38").Value = '  -3.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.67'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").Value = '0.0₆0280'
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '144.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.544'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.47%  '
$ws.Range("E51").Value = '  -2.06%  '
